# "Final tweaks at report plus Colombia"
#
# This script reproduces the meaningful content/structure changes from the
# target diff:
#   - hargeisa (sheet1): insert 3 new indicator rows (security incident,
#     ability to move freely, above-average meals) and shift the rest down.
#   - colombia_hh (sheet4): shorten 8 indicator labels.
#   - Selection / active-tab bookkeeping: hargeisa's selection moves to A8,
#     colombia_hh becomes the active tab with selection C3 (colombia_ind
#     loses the active tab but keeps its existing B24 selection).
#
# (Cosmetic-only deltas in the source diff -- defaultRowHeight/dyDescent
# font-scale values, bestFit column-width jitter, fileVersion/rupBuild,
# revisionPtr GUIDs, the author's local absPath, and raw window geometry --
# come from re-saving with a newer local Excel build and are not exposed
# through the Excel object model, so they are intentionally left alone.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. hargeisa: insert the 3 new rows.
# ---------------------------------------------------------------------
$wsHargeisa = $wb.Worksheets.Item("hargeisa")

# Insert blank rows at ascending target row numbers (2, then 5, then 7).
# Because each insert only pushes rows at/below the insertion point down,
# and every remaining target is below the current insertion point, using
# the final target row number each time lands the three blanks exactly at
# rows 2, 5 and 7 of the final 24-row sheet.
$wsHargeisa.Rows("2:2").Insert()
$wsHargeisa.Rows("5:5").Insert()
$wsHargeisa.Rows("7:7").Insert()

$wsHargeisa.Range("A2").Value = 1.1
$wsHargeisa.Range("B2").Value = "I1_sec_inc"
$wsHargeisa.Range("C2").Value = "Experience of security incident"

$wsHargeisa.Range("A5").Value = 1.2
$wsHargeisa.Range("B5").Value = "I2_move"
$wsHargeisa.Range("C5").Value = "Ability to move freely"

$wsHargeisa.Range("A7").Value = 2.1
$wsHargeisa.Range("B7").Value = "I3_meals"
$wsHargeisa.Range("C7").Value = "Above average meals per day"

# ---------------------------------------------------------------------
# 2. colombia_hh: shorten 8 indicator labels.
#
# Applied bottom-row-first to reproduce the original author's shared-
# string insertion order (new unique strings are appended to the shared
# string table in first-seen order, and the target table shows them in
# this exact row sequence: 23, 13, 21, 4, 20, 3, 24, 2).
# ---------------------------------------------------------------------
$wsColombiaHh = $wb.Worksheets.Item("colombia_hh")

$wsColombiaHh.Range("C23").Value = "Written employment contract"
$wsColombiaHh.Range("C13").Value = "Official educational establishment"
$wsColombiaHh.Range("C21").Value = "Written tenancy agreement"
$wsColombiaHh.Range("C4").Value = "Income per capita > food security line"
$wsColombiaHh.Range("C20").Value = "Income per capita > poverty line"
$wsColombiaHh.Range("C3").Value = "Natural disaster in past 12 months"
$wsColombiaHh.Range("C24").Value = "Satisfaction with current income"
$wsColombiaHh.Range("C2").Value = "Satisfaction with level of security"

# ---------------------------------------------------------------------
# 3. Selections + active tab.
#
# Select on each sheet first (Range.Select activates its sheet as a side
# effect), finishing with colombia_hh last so it ends up the active /
# tabSelected sheet, matching the workbook-level activeTab="3" target.
# colombia_ind is left untouched, so it keeps its original B24 selection
# and simply loses tabSelected once colombia_hh is activated.
# ---------------------------------------------------------------------
$wsHargeisa.Range("A8").Select()
$wsColombiaHh.Range("C3").Select()

Write-Host "edits applied"
